$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "SC-5 (2),SC-5,CM-6 b"
$ws.Range("A4").Value = "CM-5 (1),AC-6 (8),AU-12 (3),AU-7 b,AC-6 (9),AU-8 b,AU-7 a"
$ws.Range("A5").Value = "CM-7 b,AC-17 (9),AC-17 (1),CM-6 b"
$ws.Range("A8").Value = "IA-2 (12),IA-2 (11)"
$ws.Range("A10").Value = "CM-7 (2),CM-7 (5) (b)"
$ws.Range("A12").Value = "AC-7 a,AC-7 b"
$ws.Range("A15").Value = "AU-3 (1),IA-8,IA-2"
$ws.Range("A16").Value = "AC-6 (10),CM-6 b"
$ws.Range("A17").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A19").Value = "IA-5 (1) (b),IA-5 (1) (a),CM-6 b"
$ws.Range("A21").Value = "AC-12,MA-4 (7),MA-4 e,SC-10"
$ws.Range("A22").Value = "AU-3,CM-5 (1),AU-7 (1),AU-14 (1),AU-12 a,AU-6 (4),MA-4 (1) (a),CM-6 b,AU-7 a,AU-3 (1)"
$ws.Range("A25").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A29").Value = "SC-8,SC-8 (1),SC-8 (2)"
$ws.Range("A31").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AC-2 (4),AU-3 (1)"
$ws.Range("A42").Value = "SC-28 (1),SC-28"
$ws.Range("A50").Value = "IA-2 (5),CM-6 b"
$ws.Range("A53").Value = "SC-13,MA-4 (6)"
$ws.Range("A55").Value = "SC-8,AC-17 (2)"
$ws.Range("A56").Value = "MA-4 (1) (a),AU-12 c"
$ws.Range("A63").Value = "AU-5 a,AU-5 (1)"
$ws.Range("A67").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A69").Value = "CM-5 (1),AU-12 (3),AU-7 b,AU-12 a,AU-8 b,AU-12 c,CM-6 b,AU-7 a"
$ws.Range("A77").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AC-2 (4),AU-3 (1)"
$ws.Range("A80").Value = "IA-2 (4),IA-2 (3),IA-2 (2),IA-2 (1)"
$ws.Range("A81").Value = "CM-5 (3),CM-6 b"
$ws.Range("A86").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A88").Value = "CM-5 (1),AC-6 (9),AU-12 c,AC-2 (4)"
$ws.Range("A89").Value = "IA-2 (4),IA-2 (2),IA-2 (5),IA-2 (3),IA-2"
$ws.Range("A90").Value = "IA-2 (12),IA-2 (11)"
$ws.Range("A97").Value = "AU-8 b,AU-8 (1) (b),AU-8 (1) (a)"
$ws.Range("A101").Value = "IA-11,AC-3 (4)"
$ws.Range("A102").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A119").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A124").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A125").Value = "AC-18 (1),CM-7 a"
$ws.Range("A139").Value = "SI-6 b,CM-3 (5),SI-6 d"
$ws.Range("A148").Value = "AU-3,AU-14 (1),AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A157").Value = "AU-3,AU-12 a,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Range("A159").Value = "SC-8,AC-17 (2)"
$ws.Range("A181").Value = "SC-3,CM-6 b"
